# PWM et Fonction lumiere valide
# - Adds a "Points" column (H) with per-row point values
# - Marks D6 = YES, E6 = NO (the "Pilotage du pointeur lumineux" row)
# - Marks D10 = Flo (assignee for the "Commande Detection d'obstacle" row)
# - Updates the window zoom / selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data fixes on existing columns ---------------------------------------
$ws.Range("D6").Value = "YES"
$ws.Range("E6").Value = "NO"
$ws.Range("D10").Value = "Flo"

# --- New "Points" column (H) ------------------------------------------------
# Header: copy G1's formatting (centered Arial header style) onto H1, then
# set the text so it reuses the shared-string table instead of minting a
# brand-new font/style combo.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Points"

$points = [ordered]@{
  2  = 40
  3  = 30
  4  = 10
  5  = 10
  6  = 40
  7  = 10
  8  = 10
  9  = 30
  10 = 30
  11 = 10
  12 = 10
  13 = 20
  14 = 15
  15 = 25
  16 = 10
  17 = 15
  18 = 15
  19 = 20
  20 = 15
  21 = 15
  22 = 30
  23 = 40
  24 = 10
  25 = 20
  26 = 10
  27 = 10
  28 = 30
  29 = 30
}

foreach ($row in $points.Keys) {
  $ws.Range("H$row").Value = $points[$row]
}

# H6 keeps a "quote prefix" cell style (as if it had been typed as '40 and
# then corrected back to a number) - apply that formatting after the value
# is already numeric so it is not reinterpreted as text.
$ws.Range("Z1").Value = "'40"
$ws.Range("Z1").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# --- View state -------------------------------------------------------------
$excel.ActiveWindow.Zoom = 70
$ws.Range("G29").Select()

Write-Output "done"
